$d = $word.ActiveDocument

# --- Body: turn the first paragraph ("something") into a Heading 1 -------
# Apply the (localized) built-in heading style to the paragraph. Word
# recognises "Overskrift1" as the Danish id for the built-in "Heading 1"
# style and mints a style definition for it (without marking it as a
# w:customStyle, just like a genuine built-in style).
$p = $d.Paragraphs(1)
$p.Style = "Overskrift1"

# Drop the literal run text "something" - the bookmark around it stays.
[void]$d.Content.Find.Execute("something", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 2)

# --- Flesh out the minted "Overskrift1" (Heading 1) paragraph style ------
$s = $d.Styles("Overskrift1")
$s.NameLocal = "heading 1"
$s.Priority = 9
$s.NextParagraphStyle = "Normal"

$s.Font.Bold = $true
$s.Font.BoldBi = $true
$s.Font.Size = 14
$s.Font.SizeBi = 14
$s.Font.TextColor.ObjectThemeColor = 4   # wdThemeColorAccent1
$s.Font.TextColor.TintAndShade = -0.25

$s.ParagraphFormat.KeepWithNext = $true
$s.ParagraphFormat.KeepTogether = $true
$s.ParagraphFormat.SpaceBefore = 24
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.OutlineLevel = 1

# --- Companion linked character style "Overskrift 1 Tegn" ----------------
$cs = $d.Styles.Add("Overskrift1Tegn", 2)
$cs.NameLocal = "Overskrift 1 Tegn"
$cs.BaseStyle = $d.Styles("Standardskrifttypeiafsnit")
$cs.Priority = 9

$cs.Font.Bold = $true
$cs.Font.BoldBi = $true
$cs.Font.Size = 14
$cs.Font.SizeBi = 14
$cs.Font.TextColor.ObjectThemeColor = 4   # wdThemeColorAccent1
$cs.Font.TextColor.TintAndShade = -0.25

# Link the paragraph style and its character style together (both
# directions, so each one ends up with a <w:link> back to the other).
$s.LinkStyle = $cs
$cs.LinkStyle = $s
